# "Generate Report for Archive" - refresh the localization status report:
#  - flip the shared "Status" value from "Ready for handoff" to "In Translation"
#  - shrink the now-narrower Status column(s) to fit the new text

$wb = $excel.ActiveWorkbook

# 1) Replace the status text on every worksheet (Overview, zh-cn, de-de all
#    reference the same status string in their "Status"/"zh-cn"/"de-de" column).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation", 1) | Out-Null
}

# 2) Re-fit the Status columns that held that text to their new, narrower width.
$newWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = $newWidth
$wsOverview.Range("F1").ColumnWidth = $newWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = $newWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = $newWidth

$wb.Save()
